$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17-18 (Cluster tree table, first block): align style with the row above (B16/B17)
# and stamp the "<Definition>" tag in column B.
$ws.Range("B17").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B17:B18").Value = "<Definition>"

# Rows 41-59 (Cluster rule tree table, lower half): these B cells currently carry a
# stray bold/green header style; restyle them to match the plain blue style already
# used by the neighboring C column in that range, then stamp the Definition tag.
$ws.Range("C42").Copy()
$ws.Range("B41:B59").PasteSpecial(-4122)
$ws.Range("B41:B59").Value = "<Definition>"

# Rows 82-100 (Cluster rule table): column B did not exist at all; add it with the
# blue/border style used at the top of that table and stamp the Definition tag.
$ws.Range("C82").Copy()
$ws.Range("B82:B100").PasteSpecial(-4122)
$ws.Range("B82:B100").Value = "<Definition>"
